$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking columns (G:K) need to be forced to text so the stored
# cell keeps its literal text representation (e.g. "0.00", "8") instead of
# being coerced into a number by Excel.
$ws.Range("G4:K5").NumberFormat = "@"

# New row 4: duplicate of existing row 3 (Dubai (DSC) / KKR match)
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " November 01 2020"
$ws.Range("C4").Value = "KKR won by 60 runs"
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Kolkata Knight Riders"
$ws.Range("F4").Value = "Varun Aaron "
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "8"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "0.00"

# New row 5: duplicate of existing row 2 (Sharjah / Capitals match)
$ws.Range("A5").Value = " Sharjah"
$ws.Range("B5").Value = " October 09 2020"
$ws.Range("C5").Value = "Capitals won by 46 runs"
$ws.Range("D5").Value = "Rajasthan Royals"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Varun Aaron "
$ws.Range("G5").Value = "1"
$ws.Range("H5").Value = "2"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "50.00"

# Put the cell style back to the workbook default so the new rows match the
# (unstyled) formatting of the existing rows.
$ws.Range("G4:K5").Style = "Normal"
